# Updates cryptos price/volume columns (D, E) per the latest scrape.
# Numeric-looking price strings (e.g. "1.000", "9.110") are forced to stay
# literal text (matching the sheet's existing inlineStr/string cell type)
# by prefixing with an apostrophe; the cell's Style is then reset back to
# "Normal" so the quote-prefix flag picked up along the way doesn't leave
# a stray style behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = @{ V = "30.274.34"; Force = $false }; E = "  +0.24%  " },
    @{ Row = 3; D = @{ V = "1.866.54"; Force = $false }; E = "  +0.59%  " },
    @{ Row = 4; D = $null; E = "  -0.19%  " },
    @{ Row = 5; D = @{ V = "234.87"; Force = $true }; E = "  -0.29%  " },
    @{ Row = 6; D = $null; E = "  -0.28%  " },
    @{ Row = 7; D = @{ V = "0.4697"; Force = $true }; E = "  -0.17%  " },
    @{ Row = 8; D = @{ V = "0.2850"; Force = $true }; E = "  -1.46%  " },
    @{ Row = 9; D = $null; E = "  -0.02%  " },
    @{ Row = 10; D = @{ V = "21.28"; Force = $true }; E = "  -1.65%  " },
    @{ Row = 11; D = @{ V = "0.07830"; Force = $true }; E = "  -1.59%  " },
    @{ Row = 12; D = @{ V = "96.89"; Force = $true }; E = "  -0.63%  " },
    @{ Row = 13; D = @{ V = "1.871.69"; Force = $false }; E = "  +0.88%  " },
    @{ Row = 14; D = @{ V = "0.6913"; Force = $true }; E = "  +2.14%  " },
    @{ Row = 15; D = @{ V = "5.077"; Force = $true }; E = $null },
    @{ Row = 16; D = @{ V = "268.72"; Force = $true }; E = "  +0.32%  " },
    @{ Row = 17; D = @{ V = "30.266.11"; Force = $false }; E = "  +0.29%  " },
    @{ Row = 18; D = $null; E = "  +0.69%  " },
    @{ Row = 19; D = @{ V = "0.000007741"; Force = $true }; E = "  +1.97%  " },
    @{ Row = 20; D = @{ V = "1.000"; Force = $true }; E = "  -0.12%  " },
    @{ Row = 21; D = @{ V = "2.114.00"; Force = $false }; E = "  +0.80%  " },
    @{ Row = 22; D = @{ V = "1.000"; Force = $true }; E = "  -0.22%  " },
    @{ Row = 23; D = @{ V = "5.248"; Force = $true }; E = "  +0.09%  " },
    @{ Row = 24; D = @{ V = "6.149"; Force = $true }; E = "  +0.25%  " },
    @{ Row = 25; D = @{ V = "9.506"; Force = $true }; E = "  +3.96%  " },
    @{ Row = 26; D = @{ V = "165.78"; Force = $true }; E = "  -0.62%  " },
    @{ Row = 27; D = @{ V = "18.85"; Force = $true }; E = "  +0.19%  " },
    @{ Row = 28; D = @{ V = "1.934"; Force = $true }; E = "  -0.15%  " },
    @{ Row = 29; D = @{ V = "1.362"; Force = $true }; E = "  -2.78%  " },
    @{ Row = 30; D = @{ V = "0.09902"; Force = $true }; E = "  +0.37%  " },
    @{ Row = 31; D = @{ V = "4.352"; Force = $true }; E = "  +1.36%  " },
    @{ Row = 32; D = $null; E = "  -0.98%  " },
    @{ Row = 33; D = @{ V = "4.050"; Force = $true }; E = "  +1.27%  " },
    @{ Row = 34; D = @{ V = "0.04737"; Force = $true }; E = "  +1.00%  " },
    @{ Row = 35; D = $null; E = "  +0.77%  " },
    @{ Row = 36; D = @{ V = "0.7032"; Force = $true }; E = "  +0.84%  " },
    @{ Row = 37; D = $null; E = "  +0.15%  " },
    @{ Row = 38; D = @{ V = "0.01869"; Force = $true }; E = "  -0.02%  " },
    @{ Row = 39; D = @{ V = "2.773"; Force = $true }; E = "  +6.31%  " },
    @{ Row = 40; D = @{ V = "6.312"; Force = $true }; E = "  -0.24%  " },
    @{ Row = 41; D = @{ V = "73.24"; Force = $true }; E = "  -0.08%  " },
    @{ Row = 42; D = $null; E = "  +1.08%  " },
    @{ Row = 43; D = $null; E = "  -0.24%  " },
    @{ Row = 44; D = @{ V = "0.4162"; Force = $true }; E = "  +0.51%  " },
    @{ Row = 45; D = @{ V = "0.8328"; Force = $true }; E = "  -0.62%  " },
    @{ Row = 46; D = @{ V = "102.95"; Force = $true }; E = "  -0.36%  " },
    @{ Row = 47; D = @{ V = "976.74"; Force = $true }; E = "  +3.64%  " },
    @{ Row = 48; D = $null; E = "  +1.89%  " },
    @{ Row = 49; D = @{ V = "9.110"; Force = $true }; E = "  +0.19%  " },
    @{ Row = 50; D = @{ V = "34.52"; Force = $true }; E = "  +1.71%  " },
    @{ Row = 51; D = @{ V = "0.05664"; Force = $true }; E = "  +0.11%  " }
)

foreach ($item in $updates) {
    $row = $item.Row

    if ($null -ne $item.D) {
        $cell = $ws.Range("D$row")
        if ($item.D.Force) {
            $cell.Value2 = "'" + $item.D.V
            $cell.Style = "Normal"
        } else {
            $cell.Value2 = $item.D.V
        }
    }

    if ($null -ne $item.E) {
        $ws.Range("E$row").Value2 = $item.E
    }
}
